$d = $word.ActiveDocument

$pairs = @(
    @("305×3=915", "292×5=1460"),
    @("470×9=4230", "757×5=3785"),
    @("957×7=6699", "602×4=2408"),
    @("296×3=888", "220×8=1760"),
    @("770×8=6160", "549×2=1098"),
    @("324×9=2916", "883×5=4415"),
    @("628×5=3140", "929×4=3716"),
    @("231×6=1386", "478×4=1912"),
    @("814×8=6512", "625×4=2500"),
    @("780×5=3900", "734×4=2936"),
    @("310×2=620", "941×3=2823"),
    @("817×9=7353", "857×6=5142"),
    @("180×4=720", "927×9=8343"),
    @("965×8=7720", "937×8=7496"),
    @("351×8=2808", "149×3=447"),
    @("443×7=3101", "829×7=5803"),
    @("724×8=5792", "554×8=4432"),
    @("542×5=2710", "832×7=5824"),
    @("799×4=3196", "430×3=1290"),
    @("525×8=4200", "134×6=804"),
    @("758×6=4548", "445×6=2670"),
    @("622×8=4976", "595×7=4165"),
    @("349×5=1745", "539×3=1617"),
    @("488×9=4392", "396×3=1188"),
    @("444×6=2664", "308×2=616")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
